$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Row 67/68 height tweak + new empty formatted cell C67 -----------------
# Reuse the default (unstyled) style by pasting formats from a genuinely
# blank, unformatted cell so no new style entry is created.
$blank = $ws1.Range("ZZ1000")
$blank.Copy()
$ws1.Range("C67").PasteSpecial(-4122)
$ws1.Rows.Item(67).RowHeight = 13.8
$ws1.Rows.Item(68).RowHeight = 13.8

# --- New worksheet "Nova geração (09-05-18)" placed after sheet 1 ----------
$new = $wb.Worksheets.Add($null, $ws1)
$new.Name = "Nova geração (09-05-18)"

# Column widths (tuned so the engine's pixel-grid rounding lands as close as
# possible to the source workbook's 28.8 / 17.96 / 21.16 character widths)
$new.Columns.Item(1).ColumnWidth = 28.0
$new.Columns.Item(3).ColumnWidth = 17.15
$new.Columns.Item(4).ColumnWidth = 20.33

# Row heights (source sheet uses 13.8 for every populated row)
$new.Rows.Item(1).RowHeight = 13.8
$new.Rows.Item(2).RowHeight = 13.8
$new.Rows.Item(3).RowHeight = 13.8
$new.Rows.Item(4).RowHeight = 13.8
$new.Rows.Item(5).RowHeight = 13.8
$new.Rows.Item(6).RowHeight = 13.8
$new.Rows.Item(7).RowHeight = 13.8
$new.Rows.Item(8).RowHeight = 13.8
$new.Rows.Item(11).RowHeight = 13.8

# Header row
$new.Range("A1").Value = "Resource (Instance/Class)"
$new.Range("B1").Value = "#Direct Hits"
$new.Range("C1").Value = "#Indirect hits (type)"
$new.Range("D1").Value = "#Indirect hits (subclass)"
$new.Range("E1").Value = "GoodRelations (gr:) class"

# E1 reuses the bold header style already used on "Annotations per class"!C2
$ws1.Range("C2").Copy()
$new.Range("E1").PasteSpecial(-4122)

# Data rows
$new.Range("A2").Value = "Selena"
$new.Range("B2").Value = 1

$new.Range("A3").Value = "Singer"
$new.Range("B3").Value = 1
$new.Range("C3").Value = 1

$new.Range("A4").Value = "Artist"
$new.Range("C4").Value = 1

$new.Range("A5").Value = "Organization"
$new.Range("B5").Value = 1

$new.Range("A6").Value = "Person"
$new.Range("B6").Value = 1
$new.Range("C6").Value = 1

$new.Range("A7").Value = "dbpedia.org/resource/Magazine"
$new.Range("B7").Value = 3

$new.Range("A8").Value = "dbpedia.org/ontology/Magazine"
$new.Range("C8").Value = 8

$new.Range("A11").Value = "Magazine"
$new.Range("B11").Value = 3
$new.Range("C11").Value = 8

# Page setup (mirrors the source sheet's print margins / header & footer)
$ps = $new.PageSetup
$ps.LeftMargin = 56.7
$ps.RightMargin = 56.7
$ps.TopMargin = 75.8
$ps.BottomMargin = 75.8
$ps.HeaderMargin = 56.7
$ps.FooterMargin = 56.7
$ps.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ps.CenterFooter = '&"Times New Roman,Regular"&12Page &P'

# Selections / active sheet bookkeeping
$new.Range("A1:E1").Select()
$ws1.Range("A10").Select()
$new.Select()
